$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18,8).Value2 = 485.66666
$ws.Cells.Item(18,9).Value2 = 485.66666
$ws.Cells.Item(18,11).Value2 = 485.66666
$ws.Cells.Item(18,13).Value2 = -201.66666

$ws.Cells.Item(29,8).Value2 = 1056.25
$ws.Cells.Item(29,9).Value2 = 1112.5
$ws.Cells.Item(29,11).Value2 = 3337.5
$ws.Cells.Item(29,13).Value2 = -3056.5

$ws.Cells.Item(106,8).Value2 = 9462.5
$ws.Cells.Item(106,9).Value2 = 7616.6665
$ws.Cells.Item(106,11).Value2 = 7616.6665
$ws.Cells.Item(106,13).Value2 = -6985.6665

$ws.Cells.Item(112,8).Value2 = 4999.6665
$ws.Cells.Item(112,10).Value2 = 4999.6665
$ws.Cells.Item(112,12).Value2 = 14998.9995
$ws.Cells.Item(112,14).Value2 = -17214.9995

$ws.Cells.Item(118,8).Value2 = 0
$ws.Cells.Item(118,9).Value2 = 0
$ws.Cells.Item(118,11).Value2 = 0
$ws.Cells.Item(118,13).ClearContents()

$ws.Cells.Item(129,8).Value2 = 4530.75
$ws.Cells.Item(129,9).Value2 = 0
$ws.Cells.Item(129,10).Value2 = 4530.75
$ws.Cells.Item(129,11).Value2 = 0
$ws.Cells.Item(129,12).Value2 = 13592.25
$ws.Cells.Item(129,14).Value2 = -23592.25
$ws.Cells.Item(129,13).ClearContents()

$ws.Cells.Item(132,8).Value2 = 3883.7144
$ws.Cells.Item(132,9).Value2 = 3490.1538
$ws.Cells.Item(132,11).Value2 = 10470.4614
$ws.Cells.Item(132,13).Value2 = -7940.4614

$ws.Cells.Item(138,8).Value2 = 12666.667
$ws.Cells.Item(138,10).Value2 = 15000
$ws.Cells.Item(138,12).Value2 = 45000
$ws.Cells.Item(138,14).Value2 = -55280

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value2 = 918.1539
$ws.Cells.Item(32,9).Value2 = 744.6667
$ws.Cells.Item(32,11).Value2 = 744.6667
$ws.Cells.Item(32,13).Value2 = -457.6667

$ws.Cells.Item(61,8).Value2 = 7529.636
$ws.Cells.Item(61,9).Value2 = 4728.25
$ws.Cells.Item(61,11).Value2 = 4728.25
$ws.Cells.Item(61,13).Value2 = -4516.25

$ws.Cells.Item(122,8).Value2 = 5250
$ws.Cells.Item(122,10).Value2 = 5500
$ws.Cells.Item(122,12).Value2 = 16500
$ws.Cells.Item(122,14).Value2 = -21400

$ws.Cells.Item(136,8).Value2 = 7529.636
$ws.Cells.Item(136,9).Value2 = 4728.25
$ws.Cells.Item(136,11).Value2 = 14184.75
$ws.Cells.Item(136,13).Value2 = -11634.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22,8).Value2 = 321
$ws.Cells.Item(22,9).Value2 = 321
$ws.Cells.Item(22,11).Value2 = 321
$ws.Cells.Item(22,13).Value2 = -148

$ws.Cells.Item(134,8).Value2 = 11383.333
$ws.Cells.Item(134,9).Value2 = 4433.3335
$ws.Cells.Item(134,11).Value2 = 13300.0005
$ws.Cells.Item(134,13).Value2 = -10765.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value2 = 5932.231
$ws.Cells.Item(31,10).Value2 = 7282.5
$ws.Cells.Item(31,12).Value2 = 7282.5
$ws.Cells.Item(31,14).Value2 = -7872.5

$ws.Cells.Item(34,8).Value2 = 5932.231
$ws.Cells.Item(34,10).Value2 = 7282.5
$ws.Cells.Item(34,12).Value2 = 7282.5
$ws.Cells.Item(34,14).Value2 = -7686.5

$ws.Cells.Item(74,8).Value2 = 26437.666
$ws.Cells.Item(74,10).Value2 = 27156.5
$ws.Cells.Item(74,12).Value2 = 27156.5
$ws.Cells.Item(74,14).Value2 = -28904.5

$ws.Cells.Item(77,8).Value2 = 26437.666
$ws.Cells.Item(77,10).Value2 = 27156.5
$ws.Cells.Item(77,12).Value2 = 81469.5
$ws.Cells.Item(77,14).Value2 = -90205.5

$ws.Cells.Item(80,8).Value2 = 24999.666
$ws.Cells.Item(80,10).Value2 = 24999.666
$ws.Cells.Item(80,12).Value2 = 24999.666
$ws.Cells.Item(80,14).Value2 = -27245.666

$ws.Cells.Item(83,8).Value2 = 24999.666
$ws.Cells.Item(83,10).Value2 = 24999.666
$ws.Cells.Item(83,12).Value2 = 74998.99800000001
$ws.Cells.Item(83,14).Value2 = -86230.99800000001

$ws.Cells.Item(132,8).Value2 = 7452.75
$ws.Cells.Item(132,9).Value2 = 5770.3335
$ws.Cells.Item(132,11).Value2 = 17311.0005
$ws.Cells.Item(132,13).Value2 = -14781.0005

$ws.Cells.Item(134,8).Value2 = 8244.5
$ws.Cells.Item(134,9).Value2 = 3491.2
$ws.Cells.Item(134,11).Value2 = 10473.6
$ws.Cells.Item(134,13).Value2 = -7938.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38,8).Value2 = 159.33333
$ws.Cells.Item(38,9).Value2 = 127.8
$ws.Cells.Item(38,10).Value2 = 198.75
$ws.Cells.Item(38,11).Value2 = 383.4
$ws.Cells.Item(38,12).Value2 = 596.25
$ws.Cells.Item(38,13).Value2 = -36.39999999999998
$ws.Cells.Item(38,14).Value2 = -1290.25

$ws.Cells.Item(47,8).Value2 = 331.66666
$ws.Cells.Item(47,9).Value2 = 200
$ws.Cells.Item(47,10).Value2 = 397.5
$ws.Cells.Item(47,11).Value2 = 600
$ws.Cells.Item(47,12).Value2 = 1192.5
$ws.Cells.Item(47,13).Value2 = -169
$ws.Cells.Item(47,14).Value2 = -2054.5

$ws.Cells.Item(115,8).Value2 = 2000
$ws.Cells.Item(115,9).Value2 = 2000
$ws.Cells.Item(115,11).Value2 = 6000
$ws.Cells.Item(115,13).Value2 = -4825

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39,8).Value2 = 60000
$ws.Cells.Item(39,10).Value2 = 60000
$ws.Cells.Item(39,12).Value2 = 60000
$ws.Cells.Item(39,14).Value2 = -61064

$ws.Cells.Item(80,8).Value2 = 2625
$ws.Cells.Item(80,9).Value2 = 2250
$ws.Cells.Item(80,11).Value2 = 2250
$ws.Cells.Item(80,13).Value2 = -1252

$ws.Cells.Item(83,8).Value2 = 2625
$ws.Cells.Item(83,9).Value2 = 2250
$ws.Cells.Item(83,11).Value2 = 11250
$ws.Cells.Item(83,13).Value2 = -6258

$ws.Cells.Item(126,8).Value2 = 4249
$ws.Cells.Item(126,9).Value2 = 4249
$ws.Cells.Item(126,11).Value2 = 12747
$ws.Cells.Item(126,13).Value2 = -10277

$ws.Cells.Item(132,8).Value2 = 15599.4
$ws.Cells.Item(132,9).Value2 = 11332.333
$ws.Cells.Item(132,11).Value2 = 33996.999
$ws.Cells.Item(132,13).Value2 = -31466.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value2 = 2090.3333
$ws.Cells.Item(22,9).Value2 = 2657.3333
$ws.Cells.Item(22,10).Value2 = 1523.3334
$ws.Cells.Item(22,11).Value2 = 2657.3333
$ws.Cells.Item(22,12).Value2 = 1523.3334
$ws.Cells.Item(22,13).Value2 = -2362.3333
$ws.Cells.Item(22,14).Value2 = -2113.3334

$ws.Cells.Item(27,8).Value2 = 2090.3333
$ws.Cells.Item(27,9).Value2 = 2657.3333
$ws.Cells.Item(27,10).Value2 = 1523.3334
$ws.Cells.Item(27,11).Value2 = 2657.3333
$ws.Cells.Item(27,12).Value2 = 1523.3334
$ws.Cells.Item(27,13).Value2 = -2550.3333
$ws.Cells.Item(27,14).Value2 = -1737.3334

$ws.Cells.Item(46,8).Value2 = 1025
$ws.Cells.Item(46,9).Value2 = 850
$ws.Cells.Item(46,10).Value2 = 1200
$ws.Cells.Item(46,11).Value2 = 850
$ws.Cells.Item(46,12).Value2 = 1200
$ws.Cells.Item(46,13).Value2 = -662
$ws.Cells.Item(46,14).Value2 = -1576

$ws.Cells.Item(76,8).Value2 = 27333
$ws.Cells.Item(76,10).Value2 = 27333
$ws.Cells.Item(76,12).Value2 = 27333
$ws.Cells.Item(76,14).Value2 = -28009

$ws.Cells.Item(79,8).Value2 = 27333
$ws.Cells.Item(79,10).Value2 = 27333
$ws.Cells.Item(79,12).Value2 = 27333
$ws.Cells.Item(79,14).Value2 = -29673

$ws.Cells.Item(82,8).Value2 = 2714.5715
$ws.Cells.Item(82,9).Value2 = 2699.8
$ws.Cells.Item(82,10).Value2 = 2751.5
$ws.Cells.Item(82,11).Value2 = 2699.8
$ws.Cells.Item(82,12).Value2 = 2751.5
$ws.Cells.Item(82,13).Value2 = -2338.8
$ws.Cells.Item(82,14).Value2 = -3473.5

$ws.Cells.Item(85,8).Value2 = 2714.5715
$ws.Cells.Item(85,9).Value2 = 2699.8
$ws.Cells.Item(85,10).Value2 = 2751.5
$ws.Cells.Item(85,11).Value2 = 2699.8
$ws.Cells.Item(85,12).Value2 = 2751.5
$ws.Cells.Item(85,13).Value2 = -1451.8
$ws.Cells.Item(85,14).Value2 = -5247.5

$ws.Cells.Item(136,8).Value2 = 13865.571
$ws.Cells.Item(136,9).Value2 = 6388.5
$ws.Cells.Item(136,11).Value2 = 19165.5
$ws.Cells.Item(136,13).Value2 = -16615.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132,8).Value2 = 8326.076999999999
$ws.Cells.Item(132,9).Value2 = 7641.2856
$ws.Cells.Item(132,10).Value2 = 9125
$ws.Cells.Item(132,11).Value2 = 22923.8568
$ws.Cells.Item(132,12).Value2 = 27375
$ws.Cells.Item(132,13).Value2 = -20393.8568
$ws.Cells.Item(132,14).Value2 = -32435

$ws.Cells.Item(136,8).Value2 = 9790.532999999999
$ws.Cells.Item(136,9).Value2 = 8552.362999999999
$ws.Cells.Item(136,11).Value2 = 25657.089
$ws.Cells.Item(136,13).Value2 = -23107.089

Write-Host "Edits applied"